$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.081.92"
Set-TextValue $ws.Range("E2") "  +0.39%  "
Set-TextValue $ws.Range("D3") "1.888.47"
Set-TextValue $ws.Range("E3") "  -1.35%  "
Set-TextValue $ws.Range("D4") "0.9996"
Set-TextValue $ws.Range("E4") "  -0.12%  "
Set-TextValue $ws.Range("D5") "330.94"
Set-TextValue $ws.Range("E5") "  -2.35%  "
Set-TextValue $ws.Range("D6") "0.9995"
Set-TextValue $ws.Range("E6") "  -0.09%  "
Set-TextValue $ws.Range("D7") "0.4599"
Set-TextValue $ws.Range("E7") "  -2.79%  "
Set-TextValue $ws.Range("D8") "0.4072"
Set-TextValue $ws.Range("E8") "  +0.62%  "
Set-TextValue $ws.Range("D9") "47.48"
Set-TextValue $ws.Range("E9") "  -1.34%  "
Set-TextValue $ws.Range("D10") "0.07980"
Set-TextValue $ws.Range("E10") "  -2.33%  "
Set-TextValue $ws.Range("D11") "0.9934"
Set-TextValue $ws.Range("E11") "  -3.53%  "
Set-TextValue $ws.Range("D12") "21.71"
Set-TextValue $ws.Range("E12") "  -3.31%  "
Set-TextValue $ws.Range("D13") "1.882.79"
Set-TextValue $ws.Range("E13") "  -0.90%  "
Set-TextValue $ws.Range("D14") "5.913"
Set-TextValue $ws.Range("E14") "  -3.02%  "
Set-TextValue $ws.Range("D15") "7.086"
Set-TextValue $ws.Range("E15") "  -3.93%  "
Set-TextValue $ws.Range("D16") "0.9998"
Set-TextValue $ws.Range("E16") "  -0.19%  "
Set-TextValue $ws.Range("D17") "88.46"
Set-TextValue $ws.Range("E17") "  -3.36%  "
Set-TextValue $ws.Range("D18") "0.00001030"
Set-TextValue $ws.Range("E18") "  -2.26%  "
Set-TextValue $ws.Range("D19") "0.06559"
Set-TextValue $ws.Range("E19") "  -1.33%  "
Set-TextValue $ws.Range("D20") "17.47"
Set-TextValue $ws.Range("E20") "  -2.22%  "
Set-TextValue $ws.Range("D21") "0.9997"
Set-TextValue $ws.Range("E21") "  +0.00%  "
Set-TextValue $ws.Range("D22") "29.113.93"
Set-TextValue $ws.Range("E22") "  +0.43%  "
Set-TextValue $ws.Range("D23") "5.432"
Set-TextValue $ws.Range("E23") "  -2.36%  "
Set-TextValue $ws.Range("D24") "11.48"
Set-TextValue $ws.Range("E24") "  +2.49%  "
Set-TextValue $ws.Range("D25") "2.207"
Set-TextValue $ws.Range("E25") "  -2.83%  "
Set-TextValue $ws.Range("D26") "2.114.61"
Set-TextValue $ws.Range("E26") "  -0.61%  "
Set-TextValue $ws.Range("D27") "156.77"
Set-TextValue $ws.Range("E27") "  -2.54%  "
Set-TextValue $ws.Range("D28") "19.62"
Set-TextValue $ws.Range("E28") "  -2.24%  "
Set-TextValue $ws.Range("D29") "2.105"
Set-TextValue $ws.Range("E29") "  -3.50%  "
Set-TextValue $ws.Range("D30") "5.459"
Set-TextValue $ws.Range("E30") "  -1.21%  "
Set-TextValue $ws.Range("D31") "117.81"
Set-TextValue $ws.Range("E31") "  -2.62%  "
Set-TextValue $ws.Range("D32") "0.9993"
Set-TextValue $ws.Range("E32") "  -1.58%  "
Set-TextValue $ws.Range("D33") "0.09331"
Set-TextValue $ws.Range("E33") "  -2.67%  "
Set-TextValue $ws.Range("D34") "3.599"
Set-TextValue $ws.Range("E34") "  -1.70%  "
Set-TextValue $ws.Range("D35") "1.411"
Set-TextValue $ws.Range("E35") "  -1.29%  "
Set-TextValue $ws.Range("D36") "5.284"
Set-TextValue $ws.Range("E36") "  -2.62%  "
Set-TextValue $ws.Range("D37") "0.06062"
Set-TextValue $ws.Range("E37") "  -2.65%  "
Set-TextValue $ws.Range("D38") "0.02224"
Set-TextValue $ws.Range("E38") "  -2.82%  "
Set-TextValue $ws.Range("D39") "8.276"
Set-TextValue $ws.Range("E39") "  -4.54%  "
Set-TextValue $ws.Range("D40") "1.176"
Set-TextValue $ws.Range("E40") "  -2.14%  "
Set-TextValue $ws.Range("D41") "0.9986"
Set-TextValue $ws.Range("E41") "  -0.02%  "
Set-TextValue $ws.Range("D42") "0.5787"
Set-TextValue $ws.Range("E42") "  -4.23%  "
Set-TextValue $ws.Range("D43") "0.1825"
Set-TextValue $ws.Range("E43") "  -4.11%  "
Set-TextValue $ws.Range("D44") "10.13"
Set-TextValue $ws.Range("E44") "  -4.45%  "
Set-TextValue $ws.Range("D45") "1.259"
Set-TextValue $ws.Range("E45") "  -1.03%  "
Set-TextValue $ws.Range("D46") "0.07512"
Set-TextValue $ws.Range("E46") "  +2.63%  "
Set-TextValue $ws.Range("E47") "  -2.71%  "
Set-TextValue $ws.Range("E48") "  +4.95%  "
Set-TextValue $ws.Range("D49") "0.5450"
Set-TextValue $ws.Range("E49") "  -3.43%  "
Set-TextValue $ws.Range("D50") "1.903"
Set-TextValue $ws.Range("E50") "  -3.99%  "
Set-TextValue $ws.Range("D51") "111.06"
Set-TextValue $ws.Range("E51") "  -1.96%  "
